# Update for 24 April
# Adds a new "4/23/20" data column (AP) to the state-deaths table, mirroring
# the formatting of the prior date column (AO) and then overwriting the
# rows whose death toll changed from the previous day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Seed column AP with the same formatting (styles/borders/number format)
# as column AO, including the header cell, by copying AO1:AO54 -> AP1:AP54.
# This also pre-fills AP with AO's values, which are correct for every row
# that did not change day-over-day.
$ws.Range("AO1:AO54").Copy($ws.Range("AP1:AP54"))

# New header for the added date column.
$ws.Range("AP1").Value = " 4/23/20"

# Updated cumulative death counts for 4/23/20 (state rows that changed).
$newValues = @{
    4  = 249
    5  = 45
    6  = 1530
    7  = 552
    8  = 1639
    9  = 92
    10 = 139
    11 = 987
    12 = 881
    16 = 1688
    17 = 706
    18 = 96
    19 = 112
    20 = 191
    21 = 1599
    22 = 44
    23 = 748
    24 = 2360
    25 = 2977
    26 = 200
    27 = 201
    28 = 243
    30 = 47
    31 = 189
    32 = 51
    33 = 5428
    34 = 78
    35 = 20861
    36 = 281
    37 = 15
    38 = 656
    39 = 179
    40 = 83
    41 = 1685
    42 = 69
    43 = 189
    44 = 150
    46 = 170
    47 = 576
    48 = 35
    49 = 43
    50 = 372
    51 = 711
    52 = 31
    53 = 257
    54 = 7
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 42).Value = $newValues[$row]
}

# Match the author's final selection (active cell on the newly added column).
$ws.Range("AP2").Select()
